$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04941833333333333
$ws.Range("H2").Value = 0.148255
$ws.Range("I2").Value = 0.005167549122999764
$ws.Range("J2").Value = 0.005167549122999764
$ws.Range("M2").Value = 8.142376000000001
$ws.Range("N2").Value = 24.427128
$ws.Range("O2").Value = 0.1741313933276368
$ws.Range("P2").Value = 0.1741313933276368
$ws.Range("Q2").Value = 0.4023826512933334
$ws.Range("R2").Value = 3.62144386164
$ws.Range("S2").Value = 0.0008998325288769563
$ws.Range("T2").Value = 0.0008998325288769563
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04941833333333333
$ws.Range("H3").Value = 0.148255
$ws.Range("I3").Value = 0.005167549122999764
$ws.Range("J3").Value = 0.005167549122999764
$ws.Range("O3").Value = 0.5205382400466131
$ws.Range("P3").Value = 0.5205382400466131
$ws.Range("Q3").Value = 1.202859249712778
$ws.Range("R3").Value = 10.825733247415
$ws.Range("S3").Value = 0.002689906925840716
$ws.Range("T3").Value = 0.002689906925840716
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04941833333333333
$ws.Range("H4").Value = 0.148255
$ws.Range("I4").Value = 0.005167549122999764
$ws.Range("J4").Value = 0.005167549122999764
$ws.Range("O4").Value = 0.3053303666257501
$ws.Range("P4").Value = 0.3053303666257501
$ws.Range("Q4").Value = 0.70555710889
$ws.Range("R4").Value = 6.35001398001
$ws.Range("S4").Value = 0.001577809668282091
$ws.Range("T4").Value = 0.001577809668282091
$ws.Range("I5").Value = 0.806706161560336
$ws.Range("J5").Value = 0.806706161560336
$ws.Range("M5").Value = 8.142376000000001
$ws.Range("N5").Value = 24.427128
$ws.Range("O5").Value = 0.1741313933276368
$ws.Range("P5").Value = 0.1741313933276368
$ws.Range("Q5").Value = 62.81596098594668
$ws.Range("R5").Value = 565.3436488735201
$ws.Range("S5").Value = 0.140472867918491
$ws.Range("T5").Value = 0.140472867918491
$ws.Range("I6").Value = 0.806706161560336
$ws.Range("J6").Value = 0.806706161560336
$ws.Range("O6").Value = 0.5205382400466131
$ws.Range("P6").Value = 0.5205382400466131
$ws.Range("S6").Value = 0.419921405573376
$ws.Range("T6").Value = 0.419921405573376
$ws.Range("I7").Value = 0.806706161560336
$ws.Range("J7").Value = 0.806706161560336
$ws.Range("O7").Value = 0.3053303666257501
$ws.Range("P7").Value = 0.3053303666257501
$ws.Range("S7").Value = 0.246311888068469
$ws.Range("T7").Value = 0.246311888068469
$ws.Range("I8").Value = 0.1881262893166642
$ws.Range("J8").Value = 0.1881262893166643
$ws.Range("M8").Value = 8.142376000000001
$ws.Range("N8").Value = 24.427128
$ws.Range("O8").Value = 0.1741313933276368
$ws.Range("P8").Value = 0.1741313933276368
$ws.Range("Q8").Value = 14.64886995196533
$ws.Range("R8").Value = 131.839829567688
$ws.Range("S8").Value = 0.03275869288026885
$ws.Range("T8").Value = 0.03275869288026886
$ws.Range("I9").Value = 0.1881262893166642
$ws.Range("J9").Value = 0.1881262893166643
$ws.Range("O9").Value = 0.5205382400466131
$ws.Range("P9").Value = 0.5205382400466131
$ws.Range("S9").Value = 0.09792692754739635
$ws.Range("T9").Value = 0.09792692754739636
$ws.Range("I10").Value = 0.1881262893166642
$ws.Range("J10").Value = 0.1881262893166643
$ws.Range("O10").Value = 0.3053303666257501
$ws.Range("P10").Value = 0.3053303666257501
$ws.Range("S10").Value = 0.05744066888899903
$ws.Range("T10").Value = 0.05744066888899903
